# Update capital structure database values for rows 2 and 3 (earnings_debt sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("D2:D3").Value = 0.0115
$ws.Range("E2:E3").Value = 0.328
$ws.Range("G2:G3").Value = 0.08257372654155497
$ws.Range("H2:H3").Value = 0.08257372654155497
$ws.Range("I2:I3").Value = 0.08150134048257372
$ws.Range("J2:J3").Value = 0.0671587026231404
$ws.Range("K2:K3").Value = 33.6
$ws.Range("L2:L3").Value = 0.06005361930294906
$ws.Range("M2:M3").Value = 0.717
$ws.Range("N2:N3").Value = 0.003773684210526316
$ws.Range("O2:O3").Value = 0.02133928571428571
$ws.Range("P2:P3").Value = 0.717
$ws.Range("Q2:Q3").Value = 0.003773684210526316
$ws.Range("R2:R3").Value = 0.02133928571428571
$ws.Range("U2:U3").Value = 19.8
$ws.Range("V2:V3").Value = 0.1042105263157895
$ws.Range("W2:W3").Value = 0.06305122912366298
$ws.Range("X2:X3").Value = 0.1060154029432995
$ws.Range("Y2:Y3").Value = -0.04296417381963656
$ws.Range("Z2:Z3").Value = 0.9478231407758767
$ws.Range("AA2:AA3").Value = 0.06365457245069805
$ws.Range("AB2:AB3").Value = 0.08315029335186631
$ws.Range("AC2:AC3").Value = -0.01949572090116826
$ws.Range("AD2:AD3").Value = 78.3
$ws.Range("AE2:AE3").Value = 0
$ws.Range("AF2:AF3").Value = 78.3
$ws.Range("AG2:AG3").Value = 58.5
$ws.Range("AH2:AH3").Value = 0.2918374953410361
$ws.Range("AI2:AI3").Value = 0.1187263078089462
$ws.Range("AJ2:AJ3").Value = 0.2354124748490946
$ws.Range("AK2:AK3").Value = 0.09144911677348756
$ws.Range("AL2:AL3").Value = 4.03
$ws.Range("AM2:AM3").Value = 4.03
$ws.Range("AN2:AN3").Value = 1.673076923076923
$ws.Range("AO2:AO3").Value = 11.3151364764268
$ws.Range("AP2:AP3").Value = 1.25
$ws.Range("AQ2:AQ3").Value = 11.3151364764268
